$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.579703092575073
$ws.Range("B1").Value = 4.765843391418457
$ws.Range("C1").Value = 6.641824722290039
$ws.Range("D1").Value = 6.506769180297852
$ws.Range("E1").Value = 5.467338085174561
